$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content so the shared-string table can be rebuilt from scratch,
# in a deliberate order (headers, then province names top-to-bottom reflecting
# the refreshed ranking, then numeric stats, with the "last updated" banner last).
$ws.Range("A1:E64").ClearContents()

# Column headers (row 3)
$ws.Range("A3").Value = "Ciudad"
$ws.Range("B3").Value = "Casos totales"
$ws.Range("C3").Value = "Casos activos"
$ws.Range("D3").Value = "Recuperados"
$ws.Range("E3").Value = "Muertes"

# Province / city names (column A), top to bottom, in final ranked order
$ws.Range("A4").Value = "Madrid"
$ws.Range("A5").Value = "Cataluña"
$ws.Range("A6").Value = "Castilla-La Mancha"
$ws.Range("A7").Value = "Valencia/Valencia"
$ws.Range("A8").Value = "Bizkaia/Vizcaya"
$ws.Range("A9").Value = "Araba/Alava"
$ws.Range("A10").Value = "Navarra"
$ws.Range("A11").Value = "La Rioja"
$ws.Range("A12").Value = "Asturias"
$ws.Range("A13").Value = "Aragon"
$ws.Range("A14").Value = "Malaga"
$ws.Range("A15").Value = "A Coruña"
$ws.Range("A16").Value = "Alacant/Alicante"
$ws.Range("A17").Value = "Ciudad Real"
$ws.Range("A18").Value = "Toledo"
$ws.Range("A19").Value = "Pontevedra"
$ws.Range("A20").Value = "Salamanca"
$ws.Range("A21").Value = "Granada"
$ws.Range("A22").Value = "Murcia"
$ws.Range("A23").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("A24").Value = "Tenerife"
$ws.Range("A25").Value = "Albacete"
$ws.Range("A26").Value = "Sevilla"
$ws.Range("A27").Value = "Cantabria"
$ws.Range("A28").Value = "Caceres"
$ws.Range("A29").Value = "Valladolid"
$ws.Range("A30").Value = "Burgos"
$ws.Range("A31").Value = "Leon"
$ws.Range("A32").Value = "Zaragoza"
$ws.Range("A33").Value = "Segovia"
$ws.Range("A34").Value = "Guadalajara"
$ws.Range("A35").Value = "Jaen"
$ws.Range("A36").Value = "Cadiz"
$ws.Range("A37").Value = "Cordoba"
$ws.Range("A38").Value = "Castello/Castellon"
$ws.Range("A39").Value = "Badajoz"
$ws.Range("A40").Value = "Mallorca"
$ws.Range("A41").Value = "Avila"
$ws.Range("A42").Value = "Soria"
$ws.Range("A43").Value = "Gran Canaria"
$ws.Range("A44").Value = "Ourense"
$ws.Range("A45").Value = "Cuenca"
$ws.Range("A46").Value = "Zamora"
$ws.Range("A47").Value = "Lugo"
$ws.Range("A48").Value = "Almeria"
$ws.Range("A49").Value = "Palencia"
$ws.Range("A50").Value = "Huelva"
$ws.Range("A51").Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Range("A52").Value = "Teruel"
$ws.Range("A53").Value = "Huesca"
$ws.Range("A54").Value = "Melilla"
$ws.Range("A55").Value = "La Palma"
$ws.Range("A56").Value = "Ibiza"
$ws.Range("A57").Value = "Fuerteventura"
$ws.Range("A58").Value = "Menorca"
$ws.Range("A59").Value = "Lanzarote"
$ws.Range("A60").Value = "Ceuta"
$ws.Range("A61").Value = "Arroyo de la Luz"
$ws.Range("A62").Value = "El Hierro"
$ws.Range("A63").Value = "La Gomera"
$ws.Range("A64").Value = "Formentera"

# Numeric statistics (columns B-E) for every row
$ws.Range("B4").Value = 12352
$ws.Range("C4").Value = 2291
$ws.Range("D4").Value = 8526
$ws.Range("E4").Value = 1535
$ws.Range("B5").Value = 9937
$ws.Range("C5").Value = 1274
$ws.Range("D5").Value = 8147
$ws.Range("E5").Value = 516
$ws.Range("B6").Value = 2465
$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 2196
$ws.Range("E6").Value = 216
$ws.Range("B7").Value = 1317
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = 1248
$ws.Range("E7").Value = 46
$ws.Range("B8").Value = 1189
$ws.Range("C8").Value = 344
$ws.Range("D8").Value = 1145
$ws.Range("E8").Value = 44
$ws.Range("B9").Value = 1086
$ws.Range("C9").Value = 344
$ws.Range("D9").Value = 1013
$ws.Range("E9").Value = 73
$ws.Range("B10").Value = 1014
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = 972
$ws.Range("E10").Value = 31
$ws.Range("B11").Value = 802
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = 748
$ws.Range("E11").Value = 30
$ws.Range("B12").Value = 779
$ws.Range("C12").Value = 35
$ws.Range("D12").Value = 719
$ws.Range("E12").Value = 25
$ws.Range("B13").Value = 758
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 718
$ws.Range("E13").Value = 37
$ws.Range("B14").Value = 665
$ws.Range("C14").Value = 72
$ws.Range("D14").Value = 639
$ws.Range("E14").Value = 26
$ws.Range("B15").Value = 635
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 621
$ws.Range("E15").Value = 14
$ws.Range("B16").Value = 628
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 555
$ws.Range("E16").Value = 61
$ws.Range("B17").Value = 505
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 457
$ws.Range("E17").Value = 40
$ws.Range("B18").Value = 501
$ws.Range("C18").Value = 22
$ws.Range("D18").Value = 451
$ws.Range("E18").Value = 28
$ws.Range("B19").Value = 494
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 491
$ws.Range("E19").Value = 3
$ws.Range("B20").Value = 483
$ws.Range("C20").Value = 42
$ws.Range("D20").Value = 403
$ws.Range("E20").Value = 38
$ws.Range("B21").Value = 480
$ws.Range("C21").Value = 72
$ws.Range("D21").Value = 459
$ws.Range("E21").Value = 21
$ws.Range("B22").Value = 477
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 467
$ws.Range("E22").Value = 6
$ws.Range("B23").Value = 453
$ws.Range("C23").Value = 344
$ws.Range("D23").Value = 437
$ws.Range("E23").Value = 16
$ws.Range("B24").Value = 438
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 394
$ws.Range("E24").Value = 21
$ws.Range("B25").Value = 430
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 390
$ws.Range("E25").Value = 32
$ws.Range("B26").Value = 427
$ws.Range("C26").Value = 72
$ws.Range("D26").Value = 413
$ws.Range("E26").Value = 13
$ws.Range("B27").Value = 425
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 404
$ws.Range("E27").Value = 9
$ws.Range("B28").Value = 419
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 392
$ws.Range("E28").Value = 24
$ws.Range("B29").Value = 410
$ws.Range("C29").Value = 24
$ws.Range("D29").Value = 369
$ws.Range("E29").Value = 17
$ws.Range("B30").Value = 392
$ws.Range("C30").Value = 41
$ws.Range("D30").Value = 327
$ws.Range("E30").Value = 24
$ws.Range("B31").Value = 362
$ws.Range("C31").Value = 21
$ws.Range("D31").Value = 317
$ws.Range("E31").Value = 24
$ws.Range("B32").Value = 329
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 315
$ws.Range("E32").Value = 14
$ws.Range("B33").Value = 271
$ws.Range("C33").Value = 32
$ws.Range("D33").Value = 212
$ws.Range("E33").Value = 27
$ws.Range("B34").Value = 263
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 257
$ws.Range("E34").Value = 4
$ws.Range("B35").Value = 259
$ws.Range("C35").Value = 72
$ws.Range("D35").Value = 247
$ws.Range("E35").Value = 12
$ws.Range("B36").Value = 245
$ws.Range("C36").Value = 72
$ws.Range("D36").Value = 242
$ws.Range("E36").Value = 3
$ws.Range("B37").Value = 234
$ws.Range("C37").Value = 72
$ws.Range("D37").Value = 228
$ws.Range("E37").Value = 6
$ws.Range("B38").Value = 229
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 220
$ws.Range("E38").Value = 8
$ws.Range("B39").Value = 217
$ws.Range("C39").Value = 5
$ws.Range("D39").Value = 209
$ws.Range("E39").Value = 3
$ws.Range("B40").Value = 210
$ws.Range("C40").Value = 18
$ws.Range("D40").Value = 194
$ws.Range("E40").Value = 12
$ws.Range("B41").Value = 201
$ws.Range("C41").Value = 23
$ws.Range("D41").Value = 163
$ws.Range("E41").Value = 15
$ws.Range("B42").Value = 179
$ws.Range("C42").Value = 14
$ws.Range("D42").Value = 152
$ws.Range("E42").Value = 13
$ws.Range("B43").Value = 158
$ws.Range("C43").Value = 8
$ws.Range("D43").Value = 156
$ws.Range("E43").Value = 21
$ws.Range("B44").Value = 152
$ws.Range("C44").Value = 19
$ws.Range("D44").Value = 150
$ws.Range("E44").Value = 2
$ws.Range("B45").Value = 120
$ws.Range("C45").Value = 8
$ws.Range("D45").Value = 104
$ws.Range("E45").Value = 8
$ws.Range("B46").Value = 100
$ws.Range("C46").Value = 11
$ws.Range("D46").Value = 82
$ws.Range("E46").Value = 7
$ws.Range("B47").Value = 97
$ws.Range("C47").Value = 19
$ws.Range("D47").Value = 94
$ws.Range("E47").Value = 3
$ws.Range("B48").Value = 91
$ws.Range("C48").Value = 72
$ws.Range("D48").Value = 86
$ws.Range("E48").Value = 5
$ws.Range("B49").Value = 72
$ws.Range("C49").Value = 13
$ws.Range("D49").Value = 59
$ws.Range("E49").Value = 0
$ws.Range("B50").Value = 70
$ws.Range("C50").Value = 72
$ws.Range("D50").Value = 69
$ws.Range("E50").Value = 1
$ws.Range("B51").Value = 58
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 58
$ws.Range("E51").Value = 3
$ws.Range("B52").Value = 47
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 45
$ws.Range("E52").Value = 2
$ws.Range("B53").Value = 37
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 37
$ws.Range("E53").Value = 0
$ws.Range("B54").Value = 28
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 28
$ws.Range("E54").Value = 0
$ws.Range("B55").Value = 24
$ws.Range("C55").Value = 8
$ws.Range("D55").Value = 24
$ws.Range("E55").Value = 21
$ws.Range("B56").Value = 21
$ws.Range("C56").Value = 18
$ws.Range("D56").Value = 20
$ws.Range("E56").Value = 1
$ws.Range("B57").Value = 18
$ws.Range("C57").Value = 8
$ws.Range("D57").Value = 18
$ws.Range("E57").Value = 21
$ws.Range("B58").Value = 15
$ws.Range("C58").Value = 18
$ws.Range("D58").Value = 13
$ws.Range("E58").Value = 0
$ws.Range("B59").Value = 13
$ws.Range("C59").Value = 8
$ws.Range("D59").Value = 13
$ws.Range("E59").Value = 21
$ws.Range("B60").Value = 9
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 9
$ws.Range("E60").Value = 0
$ws.Range("B61").Value = 7
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 7
$ws.Range("E61").Value = 0
$ws.Range("B62").Value = 3
$ws.Range("C62").Value = 8
$ws.Range("D62").Value = 3
$ws.Range("E62").Value = 21
$ws.Range("B63").Value = 3
$ws.Range("C63").Value = 8
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 21
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 10
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 8

# "Last updated" banner (row 1) - set last so it lands at the end of the shared-string table
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 23:46"
